# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the price list (column D) for rows 31-34
$ws.Range("D31").Value = 3789.125
$ws.Range("D32").Value = 4996.397
$ws.Range("D33").Value = 6321.77
$ws.Range("D34").Value = 7128.801
